$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 571; everything from row 571 downward shifts down by one.
$ws.Rows.Item(571).Insert()

# Populate the newly inserted row 571 with the new record.
$ws.Cells.Item(571, 1).Value = 4
$ws.Cells.Item(571, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(571, 3).Value = "Los Lagos"
$ws.Cells.Item(571, 4).Value = 45223
$ws.Cells.Item(571, 5).Value = 10
$ws.Cells.Item(571, 6).Value = 100112008
$ws.Cells.Item(571, 7).Value = "Coliflor"
$ws.Cells.Item(571, 8).Value = "Sin especificar"
$ws.Cells.Item(571, 9).Value = "Primera"
$ws.Cells.Item(571, 10).Value = 1500
$ws.Cells.Item(571, 11).Value = 1500
$ws.Cells.Item(571, 12).Value = 1500
$ws.Cells.Item(571, 13).Value = 1500
$ws.Cells.Item(571, 14).Value = "$/unidad"
$ws.Cells.Item(571, 15).Value = "Región Metropolitana"
$ws.Cells.Item(571, 16).Value = 1500
$ws.Cells.Item(571, 17).Value = 1
$ws.Cells.Item(571, 18).Value = "Hortaliza"
